# Append a new match row (row 79) to the Croatia Prva NL 2023-2024 sheet,
# mirroring the existing data rows' layout/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 79
$prevRow = 78

# Copy formatting from the row above so the new row matches the existing
# styling (bold/bordered index column, date-formatted match-date column, etc.)
$ws.Range("A$prevRow`:V$prevRow").Copy() | Out-Null
$ws.Range("A$newRow`:V$newRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item($newRow, 1).Value = 78
$ws.Cells.Item($newRow, 2).Value = "croatia"
$ws.Cells.Item($newRow, 3).Value = "prva-nl"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45236.75
$ws.Cells.Item($newRow, 6).Value = "Sibenik"
$ws.Cells.Item($newRow, 7).Value = 1
$ws.Cells.Item($newRow, 8).Value = "Solin"
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 1.32
$ws.Cells.Item($newRow, 11).Value = "05/11/2023 06:12"
$ws.Cells.Item($newRow, 12).Value = 1.27
$ws.Cells.Item($newRow, 13).Value = "06/11/2023 17:52"
$ws.Cells.Item($newRow, 14).Value = 4.78
$ws.Cells.Item($newRow, 15).Value = "05/11/2023 06:12"
$ws.Cells.Item($newRow, 16).Value = 4.95
$ws.Cells.Item($newRow, 17).Value = "06/11/2023 17:53"
$ws.Cells.Item($newRow, 18).Value = 6.62
$ws.Cells.Item($newRow, 19).Value = "05/11/2023 06:12"
$ws.Cells.Item($newRow, 20).Value = 12.21
$ws.Cells.Item($newRow, 21).Value = "06/11/2023 17:56"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/croatia/prva-nl/sibenik-solin/htWRwbP8/"
